# Scheduled runner: refresh market-price columns (H:N) across all job sheets
# with current Universalis price data. Values sourced from live market API.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4255.9
$ws.Range("I32").Value = 5123.3335
$ws.Range("J32").Value = 3884.1428
$ws.Range("K32").Value = 5123.3335
$ws.Range("L32").Value = 3884.1428
$ws.Range("M32").Value = -4797.3335
$ws.Range("N32").Value = -4536.1428
$ws.Range("H43").Value = 5072.1665
$ws.Range("J43").Value = 5025
$ws.Range("L43").Value = 5025
$ws.Range("N43").Value = -5163
$ws.Range("H87").Value = 38898.945
$ws.Range("J87").Value = 39422.41
$ws.Range("L87").Value = 39422.41
$ws.Range("N87").Value = -41918.41
$ws.Range("H90").Value = 38898.945
$ws.Range("J90").Value = 39422.41
$ws.Range("L90").Value = 118267.23
$ws.Range("N90").Value = -130747.23
$ws.Range("H129").Value = 13322
$ws.Range("I129").Value = 1209.5
$ws.Range("K129").Value = 3628.5
$ws.Range("M129").Value = 1371.5
$ws.Range("H132").Value = 2538.818
$ws.Range("I132").Value = 2475.2068
$ws.Range("K132").Value = 7425.6204
$ws.Range("M132").Value = -4895.6204
$ws.Range("H137").Value = 2286.516
$ws.Range("J137").Value = 2961.6667
$ws.Range("L137").Value = 8885.000100000001
$ws.Range("N137").Value = -13985.0001
$ws.Range("H141").Value = 2326.8823
$ws.Range("I141").Value = 2326.8823
$ws.Range("K141").Value = 6980.646900000001
$ws.Range("M141").Value = -1800.646900000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2086.182
$ws.Range("I2").Value = 1859.0667
$ws.Range("J2").Value = 2572.8572
$ws.Range("K2").Value = 1859.0667
$ws.Range("L2").Value = 2572.8572
$ws.Range("M2").Value = -1746.0667
$ws.Range("N2").Value = -2798.8572
$ws.Range("H4").Value = 356.16666
$ws.Range("J4").Value = 196.66667
$ws.Range("L4").Value = 196.66667
$ws.Range("N4").Value = -428.66667
$ws.Range("H5").Value = 1540
$ws.Range("I5").Value = 1550
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 1550
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -1438
$ws.Range("N5").Value = -1724
$ws.Range("H32").Value = 1709.674
$ws.Range("I32").Value = 1503.1951
$ws.Range("K32").Value = 1503.1951
$ws.Range("M32").Value = -1216.1951
$ws.Range("H45").Value = 1399.5
$ws.Range("I45").Value = 1399.5
$ws.Range("K45").Value = 1399.5
$ws.Range("M45").Value = -1022.5
$ws.Range("H97").Value = 518.6667
$ws.Range("I97").Value = 504
$ws.Range("J97").Value = 548
$ws.Range("K97").Value = 504
$ws.Range("L97").Value = 548
$ws.Range("M97").Value = -8
$ws.Range("N97").Value = -1540
$ws.Range("H116").Value = 2086.182
$ws.Range("I116").Value = 1859.0667
$ws.Range("J116").Value = 2572.8572
$ws.Range("K116").Value = 1859.0667
$ws.Range("L116").Value = 2572.8572
$ws.Range("M116").Value = 434.9332999999999
$ws.Range("N116").Value = -7160.8572
$ws.Range("H132").Value = 5137681.5
$ws.Range("J132").Value = 11251738
$ws.Range("L132").Value = 33755214
$ws.Range("N132").Value = -33760274
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2086.182
$ws.Range("I3").Value = 1859.0667
$ws.Range("J3").Value = 2572.8572
$ws.Range("K3").Value = 1859.0667
$ws.Range("L3").Value = 2572.8572
$ws.Range("M3").Value = -1745.0667
$ws.Range("N3").Value = -2800.8572
$ws.Range("H4").Value = 1540
$ws.Range("I4").Value = 1550
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 1550
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -1435
$ws.Range("N4").Value = -1730
$ws.Range("H22").Value = 3393.125
$ws.Range("I22").Value = 4114.143
$ws.Range("K22").Value = 4114.143
$ws.Range("M22").Value = -3941.143
$ws.Range("H86").Value = 2565.2222
$ws.Range("I86").Value = 2565.2222
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2565.2222
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1442.2222
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2565.2222
$ws.Range("I89").Value = 2565.2222
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12826.111
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -7210.111000000001
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 4971.826
$ws.Range("I94").Value = 6102.4375
$ws.Range("J94").Value = 2387.5715
$ws.Range("K94").Value = 6102.4375
$ws.Range("L94").Value = 2387.5715
$ws.Range("M94").Value = -5651.4375
$ws.Range("N94").Value = -3289.5715
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 309.83334
$ws.Range("I7").Value = 496.66666
$ws.Range("J7").Value = 123
$ws.Range("K7").Value = 496.66666
$ws.Range("L7").Value = 123
$ws.Range("M7").Value = -383.66666
$ws.Range("N7").Value = -349
$ws.Range("H22").Value = 21199.8
$ws.Range("I22").Value = 33999.668
$ws.Range("K22").Value = 33999.668
$ws.Range("M22").Value = -33649.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5001578
$ws.Range("I4").Value = 7667547
$ws.Range("J4").Value = 2886.875
$ws.Range("K4").Value = 23002641
$ws.Range("L4").Value = 8660.625
$ws.Range("M4").Value = -23002529
$ws.Range("N4").Value = -8884.625
$ws.Range("H113").Value = 125594
$ws.Range("I113").Value = 250319
$ws.Range("K113").Value = 750957
$ws.Range("M113").Value = -748787
$ws.Range("H121").Value = 240258.8
$ws.Range("I121").Value = 333433
$ws.Range("J121").Value = 100497.5
$ws.Range("K121").Value = 1000299
$ws.Range("L121").Value = 301492.5
$ws.Range("M121").Value = -998989
$ws.Range("N121").Value = -304112.5
$ws.Range("H131").Value = 1530.2142
$ws.Range("I131").Value = 1043.9412
$ws.Range("K131").Value = 3131.8236
$ws.Range("M131").Value = 1908.1764
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 212.18182
$ws.Range("J2").Value = 289.8
$ws.Range("L2").Value = 289.8
$ws.Range("N2").Value = -515.8
$ws.Range("H32").Value = 56277.25
$ws.Range("J32").Value = 56277.25
$ws.Range("L32").Value = 56277.25
$ws.Range("N32").Value = -56869.25
$ws.Range("H101").Value = 13663.5
$ws.Range("J101").Value = 13663.5
$ws.Range("L101").Value = 13663.5
$ws.Range("N101").Value = -20153.5
$ws.Range("H102").Value = 4974.407
$ws.Range("I102").Value = 3321.3333
$ws.Range("J102").Value = 18199
$ws.Range("K102").Value = 3321.3333
$ws.Range("L102").Value = 18199
$ws.Range("M102").Value = -1699.3333
$ws.Range("N102").Value = -21443
$ws.Range("H105").Value = 39239.89
$ws.Range("J105").Value = 39239.89
$ws.Range("L105").Value = 39239.89
$ws.Range("N105").Value = -46227.89
$ws.Range("H132").Value = 1587572.5
$ws.Range("J132").Value = 1978.6875
$ws.Range("L132").Value = 5936.0625
$ws.Range("N132").Value = -10996.0625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2381.889
$ws.Range("J16").Value = 4409.091
$ws.Range("L16").Value = 4409.091
$ws.Range("N16").Value = -4749.091
$ws.Range("H93").Value = 2548.2856
$ws.Range("J93").Value = 3212.5
$ws.Range("L93").Value = 3212.5
$ws.Range("N93").Value = -5708.5
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H132").Value = 17867244
$ws.Range("I132").Value = 19241532
$ws.Range("K132").Value = 57724596
$ws.Range("M132").Value = -57722066
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 45900
$ws.Range("I64").Value = 45900
$ws.Range("K64").Value = 45900
$ws.Range("M64").Value = -45652
$ws.Range("H67").Value = 45900
$ws.Range("I67").Value = 45900
$ws.Range("K67").Value = 45900
$ws.Range("M67").Value = -45042
$ws.Range("H96").Value = 1911.4445
$ws.Range("I96").Value = 1516.6666
$ws.Range("J96").Value = 2701
$ws.Range("K96").Value = 1516.6666
$ws.Range("L96").Value = 2701
$ws.Range("M96").Value = -143.6666
$ws.Range("N96").Value = -5447
$ws.Range("H132").Value = 10420377
$ws.Range("I132").Value = 13514883
$ws.Range("J132").Value = 11584.728
$ws.Range("K132").Value = 40544649
$ws.Range("L132").Value = 34754.18399999999
$ws.Range("M132").Value = -40542119
$ws.Range("N132").Value = -39814.18399999999
